$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1043256666666667
$ws.Range("H2").Value = 0.312977
$ws.Range("I2").Value = 0.02547563162231953
$ws.Range("J2").Value = 0.02547563162231953
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.024701
$ws.Range("N2").Value = 9.074103000000001
$ws.Range("O2").Value = 0.1596375877334842
$ws.Range("P2").Value = 0.1596375877334843
$ws.Range("Q2").Value = 0.3155539482923334
$ws.Range("R2").Value = 2.839985534631
$ws.Range("S2").Value = 0.00406686837817396
$ws.Range("T2").Value = 0.00406686837817396

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1043256666666667
$ws.Range("H3").Value = 0.312977
$ws.Range("I3").Value = 0.02547563162231953
$ws.Range("J3").Value = 0.02547563162231953
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("O3").Value = 0.6072559333217162
$ws.Range("P3").Value = 0.6072559333217163
$ws.Range("Q3").Value = 1.200356445522889
$ws.Range("R3").Value = 10.803208009706
$ws.Range("S3").Value = 0.01547022845777188
$ws.Range("T3").Value = 0.01547022845777188

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1043256666666667
$ws.Range("H4").Value = 0.312977
$ws.Range("I4").Value = 0.02547563162231953
$ws.Range("J4").Value = 0.02547563162231953
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.368554666666666
$ws.Range("N4").Value = 13.105664
$ws.Range("O4").Value = 0.2305634602787257
$ws.Range("P4").Value = 0.2305634602787257
$ws.Range("Q4").Value = 0.4557523779697777
$ws.Range("R4").Value = 4.101771401728
$ws.Range("S4").Value = 0.005873749779628117
$ws.Range("T4").Value = 0.005873749779628118

$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1043256666666667
$ws.Range("H5").Value = 0.312977
$ws.Range("I5").Value = 0.02547563162231953
$ws.Range("J5").Value = 0.02547563162231953
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.04818333333333333
$ws.Range("N5").Value = 0.14455
$ws.Range("O5").Value = 0.002543018666073676
$ws.Range("P5").Value = 0.002543018666073677
$ws.Range("Q5").Value = 0.005026758372222222
$ws.Range("R5").Value = 0.04524082535000001
$ws.Range("S5").Value = 0.00006478500674557538
$ws.Range("T5").Value = 0.0000647850067455754

$ws.Range("I6").Value = 0.9745243683776804
$ws.Range("J6").Value = 0.9745243683776804
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.024701
$ws.Range("N6").Value = 9.074103000000001
$ws.Range("O6").Value = 0.1596375877334842
$ws.Range("P6").Value = 0.1596375877334843
$ws.Range("Q6").Value = 12.07094751202367
$ws.Range("R6").Value = 108.638527608213
$ws.Range("S6").Value = 0.1555707193553103
$ws.Range("T6").Value = 0.1555707193553103

$ws.Range("I7").Value = 0.9745243683776804
$ws.Range("J7").Value = 0.9745243683776804
$ws.Range("O7").Value = 0.6072559333217162
$ws.Range("P7").Value = 0.6072559333217163
$ws.Range("Q7").Value = 45.91747220415977
$ws.Range("R7").Value = 413.257249837438
$ws.Range("S7").Value = 0.5917857048639443
$ws.Range("T7").Value = 0.5917857048639444

$ws.Range("I8").Value = 0.9745243683776804
$ws.Range("J8").Value = 0.9745243683776804
$ws.Range("M8").Value = 4.368554666666666
$ws.Range("N8").Value = 13.105664
$ws.Range("O8").Value = 0.2305634602787257
$ws.Range("P8").Value = 0.2305634602787257
$ws.Range("Q8").Value = 17.43398573437155
$ws.Range("R8").Value = 156.905871609344
$ws.Range("S8").Value = 0.2246897104990976
$ws.Range("T8").Value = 0.2246897104990976

$ws.Range("I9").Value = 0.9745243683776804
$ws.Range("J9").Value = 0.9745243683776804
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.04818333333333333
$ws.Range("N9").Value = 0.14455
$ws.Range("O9").Value = 0.002543018666073676
$ws.Range("P9").Value = 0.002543018666073677
$ws.Range("Q9").Value = 0.1922895808944444
$ws.Range("R9").Value = 1.73060622805
$ws.Range("S9").Value = 0.002478233659328101
$ws.Range("T9").Value = 0.002478233659328101
